# "checklist new version is added"
#
# 1. login sheet: move the saved selection from G3 to D5
# 2. ChecklistManagement sheet: update the view (drop topLeftCell, move
#    selection from Q18 to G19) and fill in row 9 with the data for the
#    new "checklist new version" scenario
# 3. Leave ChecklistManagement as the active sheet/tab, as it was before.

$wb = $excel.ActiveWorkbook

# --- login sheet: saved selection moves from G3 to D5 -----------------
$login = $wb.Sheets.Item("login")
$login.Range("D5").Select()

# --- ChecklistManagement sheet -----------------------------------------
$cm = $wb.Sheets.Item("ChecklistManagement")
$cm.Select()

# Row 9 values (new checklist-new-version test case)
$cm.Range("D9").Value = "2"
$cm.Range("E9").Value = "positive"
$cm.Range("F9").Value = "checklist-management"
$cm.Range("G9").Value = "USER MANAGEMENT"
$cm.Range("H9").Value = ""
$cm.Range("I9").Value = "LOGIN"
$cm.Range("J9").Value = "1"
$cm.Range("K9").Value = "saleHeader"
$cm.Range("L9").Value = "positive"
$cm.Range("M9").Value = "1"
$cm.Range("N9").Value = "1"
$cm.Range("O9").Value = "2"
$cm.Range("P9").Value = "2"
$cm.Range("Q9").Value = "1"
$cm.Range("R9").Value = "The checklist cannot be empty"
$cm.Range("T9").Value = ""
$cm.Range("U9").Value = ""

# Writing a value through COM silently drops the cell's "quote prefix"
# formatting (the style used for text that looks numeric). Re-apply the
# original (non quote-prefixed) look of the row by pasting the format
# from neighboring cells that already carry it, without touching values.
$cm.Range("C9").Copy()
$cm.Range("D9:R9").PasteSpecial(-4122)
$cm.Range("T9:U9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Saved selection on this sheet moves from Q18 to G19 (and the
# topLeftCell freeze/scroll anchor goes away along with it).
$cm.Range("G19").Select()
